$wb = $excel.ActiveWorkbook

# --- Sheet 1: LoginTestData ------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new sheet "DashboardPageTabs" right after LoginTestData -------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "DashboardPageTabs"

# Reuse the existing header style (bold font + orange fill + border) from
# A1 on LoginTestData for the new header cell, then just recolor the fill.
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)

# Reuse the existing plain bordered body style (from B2) for the data rows.
$ws1.Range("B2").Copy()
$ws2.Range("A2:A17").PasteSpecial(-4122)

# Populate values in the same order the strings were originally authored so
# that the shared-string table ends up in the expected order: the main list
# top-to-bottom first, with Opportunities / List Emails / Notes (rows that
# were inserted afterwards) filled in last.
$ws2.Range("A1").Value = "dashboardPageTabs"
$ws2.Range("A2").Value = "Home"
$ws2.Range("A3").Value = "Accounts"
$ws2.Range("A4").Value = "Contacts"
$ws2.Range("A5").Value = "Leads"
$ws2.Range("A7").Value = "Tasks"
$ws2.Range("A8").Value = "Calendar"
$ws2.Range("A9").Value = "Dashboards"
$ws2.Range("A11").Value = "Reports"
$ws2.Range("A12").Value = "Groups"
$ws2.Range("A13").Value = "Forecasts"
$ws2.Range("A14").Value = "Files"
$ws2.Range("A16").Value = "Quotes"
$ws2.Range("A17").Value = "Chatter"
$ws2.Range("A6").Value = "Opportunities"
$ws2.Range("A15").Value = "List Emails"
$ws2.Range("A10").Value = "Notes"

# Header cell gets a distinct green fill.
$ws2.Range("A1").Interior.Color = 5287936

# Column sizing to match the authored sheet.
$ws2.Range("A1").ColumnWidth = 16.6

# Selection state left on the new sheet.
$ws2.Range("A10").Select()

# --- Back on sheet 1: drop the two stray / duplicate rows ------------------
$ws1.Rows("8:9").Delete()
$ws1.Range("C9:C10").Select()

# The new sheet becomes the active tab, as in the authored workbook.
$ws2.Activate()
